# adding fixed operating cost GEOMETRIC trend
$wb = $excel.ActiveWorkbook

$wsImport = $wb.Worksheets.Item("Import Priorities")
$wsCoupling = $wb.Worksheets.Item("Coupling Parameters")

# Insert two new rows right after the "FuelPriceTrends" row (row 2),
# pushing the remaining priority rows down by two.
$wsImport.Rows.Item(3).Resize(2).Insert()

# Bump FuelPriceTrends' own priority up to make room for the new trends.
$wsImport.Range("B2").Value = 12

# New rows: GeometricTrends (priority 11) and StepTrends (priority 10).
$wsImport.Range("A3").Value = "GeometricTrends"
$wsImport.Range("B3").Value = 11
$wsImport.Range("A4").Value = "StepTrends"
$wsImport.Range("B4").Value = 10

# Restore selection/active-sheet state recorded in the workbook views.
$wsCoupling.Range("B9").Select()
$wsImport.Activate()
$wsImport.Range("E9").Select()
